# Add a header row to the "Links" worksheet:
#   Location | Notes | Link
# This pushes all existing data down by one row, so the worksheet's
# hyperlinks (which this runtime does not auto-shift on row insert) are
# captured beforehand and re-created one row lower, preserving their
# original cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remember where each hyperlink currently lives (row only - they are
#        all in column C) before we disturb the layout.
$oldHyperlinkRows = @()
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $oldHyperlinkRows += $ws.Hyperlinks.Item($i).Range.Row
}

# --- 2. Insert a brand-new row 1 and push everything else down.
$ws.Rows("1:1").Insert()

# --- 3. Write the new header row.
$ws.Range("A1").Value2 = "Location"
$ws.Range("B1").Value2 = "Notes"
$ws.Range("C1").Value2 = "Link"

# --- 4. Re-create the hyperlinks one row below where they used to be
#        (their cell text already equals the hyperlink target, so reuse it),
#        preserving each cell's original style.
$ws.Hyperlinks.Delete()
foreach ($r in $oldHyperlinkRows) {
    $newRow = $r + 1
    $cell = $ws.Cells.Item($newRow, 3)
    $target = $cell.Value2
    $savedStyle = $cell.Style
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
    $cell.Style = $savedStyle
}

# --- 5. Match the saved selection/scroll state from the edit (cursor on C1).
$ws.Range("C1").Select()
